$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-08-25 01:01:20"

$wsZhCn.Range("H2").Value = "2016-08-25 01:01:14"
$wsZhCn.Range("K2").Value = "2016-08-25 01:01:50"

$wsDeDe.Range("K2").Value = "2016-08-25 01:01:57"
